# Apply cell value updates per the target diff (Lamia_Profits workbook)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 550
$ws.Cells.Item(41, 9).Value = 360.0909
$ws.Cells.Item(41, 11).Value = 360.0909
$ws.Cells.Item(41, 13).Value = 79.90910000000002
$ws.Cells.Item(55, 8).Value = 853.9286
$ws.Cells.Item(55, 9).Value = 380.125
$ws.Cells.Item(55, 10).Value = 1145.5
$ws.Cells.Item(55, 11).Value = 380.125
$ws.Cells.Item(55, 12).Value = 1145.5
$ws.Cells.Item(55, 13).Value = -166.125
$ws.Cells.Item(55, 14).Value = -1573.5
$ws.Cells.Item(69, 8).Value = 8766.25
$ws.Cells.Item(69, 9).Value = 8675
$ws.Cells.Item(69, 11).Value = 26025
$ws.Cells.Item(69, 13).Value = -25151
$ws.Cells.Item(72, 8).Value = 8766.25
$ws.Cells.Item(72, 9).Value = 8675
$ws.Cells.Item(72, 11).Value = 78075
$ws.Cells.Item(72, 13).Value = -73707
$ws.Cells.Item(76, 8).Value = 6408.15
$ws.Cells.Item(76, 9).Value = 5674.625
$ws.Cells.Item(76, 10).Value = 6897.1665
$ws.Cells.Item(76, 11).Value = 5674.625
$ws.Cells.Item(76, 12).Value = 6897.1665
$ws.Cells.Item(76, 13).Value = -5359.625
$ws.Cells.Item(76, 14).Value = -7527.1665
$ws.Cells.Item(79, 8).Value = 6408.15
$ws.Cells.Item(79, 9).Value = 5674.625
$ws.Cells.Item(79, 10).Value = 6897.1665
$ws.Cells.Item(79, 11).Value = 5674.625
$ws.Cells.Item(79, 12).Value = 6897.1665
$ws.Cells.Item(79, 13).Value = -4582.625
$ws.Cells.Item(79, 14).Value = -9081.166499999999
$ws.Cells.Item(96, 8).Value = 773.1875
$ws.Cells.Item(96, 10).Value = 385.5
$ws.Cells.Item(96, 12).Value = 1156.5
$ws.Cells.Item(96, 14).Value = -3902.5
$ws.Cells.Item(107, 8).Value = 808.3333
$ws.Cells.Item(107, 9).Value = 808.3333
$ws.Cells.Item(107, 11).Value = 808.3333
$ws.Cells.Item(107, 13).Value = 1111.6667
$ws.Cells.Item(112, 8).Value = 2310.6667
$ws.Cells.Item(112, 10).Value = 2310.6667
$ws.Cells.Item(112, 12).Value = 6932.000100000001
$ws.Cells.Item(112, 14).Value = -9148.000100000001
$ws.Cells.Item(137, 8).Value = 13336477
$ws.Cells.Item(137, 9).Value = 90912136
$ws.Cells.Item(137, 11).Value = 272736408
$ws.Cells.Item(137, 13).Value = -272733858
$ws.Cells.Item(138, 8).Value = 4508.7383
$ws.Cells.Item(138, 9).Value = 2652.75
$ws.Cells.Item(138, 10).Value = 5251.1333
$ws.Cells.Item(138, 11).Value = 7958.25
$ws.Cells.Item(138, 12).Value = 15753.3999
$ws.Cells.Item(138, 13).Value = -2818.25
$ws.Cells.Item(138, 14).Value = -26033.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 3292.1365
$ws.Cells.Item(74, 9).Value = 2110.25
$ws.Cells.Item(74, 11).Value = 2110.25
$ws.Cells.Item(74, 13).Value = -1236.25
$ws.Cells.Item(77, 8).Value = 3292.1365
$ws.Cells.Item(77, 9).Value = 2110.25
$ws.Cells.Item(77, 11).Value = 10551.25
$ws.Cells.Item(77, 13).Value = -6183.25
$ws.Cells.Item(122, 8).Value = 3530.5862
$ws.Cells.Item(122, 9).Value = 2988.5557
$ws.Cells.Item(122, 10).Value = 4417.5454
$ws.Cells.Item(122, 11).Value = 8965.667099999999
$ws.Cells.Item(122, 12).Value = 13252.6362
$ws.Cells.Item(122, 13).Value = -6515.667099999999
$ws.Cells.Item(122, 14).Value = -18152.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1871.3158
$ws.Cells.Item(16, 9).Value = 1152.9166
$ws.Cells.Item(16, 10).Value = 3102.8572
$ws.Cells.Item(16, 11).Value = 1152.9166
$ws.Cells.Item(16, 12).Value = 3102.8572
$ws.Cells.Item(16, 13).Value = -865.9166
$ws.Cells.Item(16, 14).Value = -3676.8572
$ws.Cells.Item(51, 8).Value = 39999.918
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 13).ClearContents()  # was -39263
$ws.Cells.Item(58, 8).Value = 4102.8076
$ws.Cells.Item(58, 9).Value = 1308.75
$ws.Cells.Item(58, 10).Value = 5344.6113
$ws.Cells.Item(58, 11).Value = 1308.75
$ws.Cells.Item(58, 12).Value = 5344.6113
$ws.Cells.Item(58, 13).Value = -1105.75
$ws.Cells.Item(58, 14).Value = -5750.6113
$ws.Cells.Item(61, 8).Value = 39999.918
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 13).ClearContents()  # was -39651
$ws.Cells.Item(113, 8).Value = 1871.3158
$ws.Cells.Item(113, 9).Value = 1152.9166
$ws.Cells.Item(113, 10).Value = 3102.8572
$ws.Cells.Item(113, 11).Value = 1152.9166
$ws.Cells.Item(113, 12).Value = 3102.8572
$ws.Cells.Item(113, 13).Value = 1017.0834
$ws.Cells.Item(113, 14).Value = -7442.8572
$ws.Cells.Item(132, 8).Value = 4338.1113
$ws.Cells.Item(132, 9).Value = 3927.5
$ws.Cells.Item(132, 11).Value = 11782.5
$ws.Cells.Item(132, 13).Value = -9252.5
$ws.Cells.Item(133, 8).Value = 55916.07
$ws.Cells.Item(133, 10).Value = 56165.91
$ws.Cells.Item(133, 12).Value = 56165.91
$ws.Cells.Item(133, 14).Value = -61225.91
$ws.Cells.Item(135, 8).Value = 59508.5
$ws.Cells.Item(135, 9).Value = 59490
$ws.Cells.Item(135, 10).Value = 59510.184
$ws.Cells.Item(135, 11).Value = 59490
$ws.Cells.Item(135, 12).Value = 59510.184
$ws.Cells.Item(135, 13).Value = -54420
$ws.Cells.Item(135, 14).Value = -69650.18400000001
$ws.Cells.Item(136, 8).Value = 4102.8076
$ws.Cells.Item(136, 9).Value = 1308.75
$ws.Cells.Item(136, 10).Value = 5344.6113
$ws.Cells.Item(136, 11).Value = 3926.25
$ws.Cells.Item(136, 12).Value = 16033.8339
$ws.Cells.Item(136, 13).Value = -1376.25
$ws.Cells.Item(136, 14).Value = -21133.8339

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 2366.2903
$ws.Cells.Item(122, 10).Value = 3172.7144
$ws.Cells.Item(122, 12).Value = 28554.4296
$ws.Cells.Item(122, 14).Value = -33454.4296

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2851.3635
$ws.Cells.Item(102, 9).Value = 1065
$ws.Cells.Item(102, 11).Value = 1065
$ws.Cells.Item(102, 13).Value = 557
$ws.Cells.Item(126, 8).Value = 3665.9656
$ws.Cells.Item(126, 9).Value = 2472.1667
$ws.Cells.Item(126, 10).Value = 5619.4546
$ws.Cells.Item(126, 11).Value = 7416.500100000001
$ws.Cells.Item(126, 12).Value = 16858.3638
$ws.Cells.Item(126, 13).Value = -4946.500100000001
$ws.Cells.Item(126, 14).Value = -21798.3638
$ws.Cells.Item(132, 8).Value = 6033.6875
$ws.Cells.Item(132, 9).Value = 1680.5454
$ws.Cells.Item(132, 11).Value = 5041.6362
$ws.Cells.Item(132, 13).Value = -2511.6362
$ws.Cells.Item(134, 8).Value = 65518.75
$ws.Cells.Item(134, 10).Value = 65518.75
$ws.Cells.Item(134, 12).Value = 196556.25
$ws.Cells.Item(134, 14).Value = -201626.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7412.815
$ws.Cells.Item(7, 9).Value = 5944.8096
$ws.Cells.Item(7, 11).Value = 5944.8096
$ws.Cells.Item(7, 13).Value = -5832.8096
$ws.Cells.Item(55, 8).Value = 3203.963
$ws.Cells.Item(55, 10).Value = 4486
$ws.Cells.Item(55, 12).Value = 4486
$ws.Cells.Item(55, 14).Value = -4832
$ws.Cells.Item(82, 8).Value = 4640.4062
$ws.Cells.Item(82, 9).Value = 2272.5
$ws.Cells.Item(82, 11).Value = 2272.5
$ws.Cells.Item(82, 13).Value = -1911.5
$ws.Cells.Item(85, 8).Value = 4640.4062
$ws.Cells.Item(85, 9).Value = 2272.5
$ws.Cells.Item(85, 11).Value = 2272.5
$ws.Cells.Item(85, 13).Value = -1024.5
$ws.Cells.Item(100, 8).Value = 5606.4
$ws.Cells.Item(100, 9).Value = 4521.6
$ws.Cells.Item(100, 11).Value = 4521.6
$ws.Cells.Item(100, 13).Value = -3980.6
$ws.Cells.Item(126, 8).Value = 7412.815
$ws.Cells.Item(126, 9).Value = 5944.8096
$ws.Cells.Item(126, 11).Value = 17834.4288
$ws.Cells.Item(126, 13).Value = -15364.4288

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 6921.4287
$ws.Cells.Item(62, 10).Value = 7312.5
$ws.Cells.Item(62, 12).Value = 7312.5
$ws.Cells.Item(62, 14).Value = -8560.5
$ws.Cells.Item(65, 8).Value = 6921.4287
$ws.Cells.Item(65, 10).Value = 7312.5
$ws.Cells.Item(65, 12).Value = 36562.5
$ws.Cells.Item(65, 14).Value = -42802.5
$ws.Cells.Item(81, 8).Value = 6612.5
$ws.Cells.Item(81, 9).Value = 4842.5713
$ws.Cells.Item(81, 11).Value = 9685.142599999999
$ws.Cells.Item(81, 13).Value = -8624.142599999999
$ws.Cells.Item(84, 8).Value = 6612.5
$ws.Cells.Item(84, 9).Value = 4842.5713
$ws.Cells.Item(84, 11).Value = 48425.713
$ws.Cells.Item(84, 13).Value = -43121.713
$ws.Cells.Item(113, 8).Value = 251.32
$ws.Cells.Item(113, 9).Value = 264.47827
$ws.Cells.Item(113, 11).Value = 793.43481
$ws.Cells.Item(113, 13).Value = 1376.56519
$ws.Cells.Item(138, 8).Value = 75429
$ws.Cells.Item(138, 10).Value = 75429
$ws.Cells.Item(138, 12).Value = 75429
$ws.Cells.Item(138, 14).Value = -85709

Write-Host "Applied all cell updates."
